$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 2 (existing incident row) - update in place
# ------------------------------------------------------------------

# B2 now holds a single "Date de creation" value; C2 ("Date de cloture") is removed entirely.
$ws.Range("B2").Value = 45819.41148925926
$ws.Range("C2").Clear()

# NumRef (A2) must stay text so the leading zero survives.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "06250003"

$ws.Range("D2").Value = "N/C"
$ws.Range("E2").Value = "Type d'incident 1"
$ws.Range("F2").Value = ""

# Equipement (G2) is no longer populated.
$ws.Range("G2").Clear()

$ws.Range("H2").Value = "P19"
$ws.Range("I2").Value = "Shift(14-22)"
$ws.Range("J2").Value = "Admin User"
$ws.Range("K2").Value = "--"
$ws.Range("L2").Value = "--"
$ws.Range("M2").Value = "This incident happened due to negligence of the Guerit chef"
$ws.Range("O2").Value = "EN MAINTENANCE"

# ------------------------------------------------------------------
# Row 3 (new incident row)
# ------------------------------------------------------------------

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "06250004"

$ws.Range("B3").NumberFormat = "m/d/yy"
$ws.Range("B3").Value = 45820.635023090275
$ws.Range("C3").NumberFormat = "m/d/yy"
$ws.Range("C3").Value = 45820.63515135417

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "0"

$ws.Range("E3").Value = "Type d'incident 1"
$ws.Range("F3").Value = "Cause incident 1"
$ws.Range("H3").Value = "P11"
$ws.Range("I3").Value = "Shift(14-22)"
$ws.Range("J3").Value = "Admin User"
$ws.Range("K3").Value = "NOELLE JEANNE"
$ws.Range("L3").Value = "Admin User"
$ws.Range("M3").Value = ""
$ws.Range("O3").Value = "CLOTURE"

# ------------------------------------------------------------------
# Row 4 (new incident row)
# ------------------------------------------------------------------

$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "06250001"

$ws.Range("B4").NumberFormat = "m/d/yy"
$ws.Range("B4").Value = 45818.687470613426

$ws.Range("D4").Value = "N/C"
$ws.Range("E4").Value = "Type d'incident 1"
$ws.Range("F4").Value = ""
$ws.Range("H4").Value = "SALLE D'ATTENTE NIVEAU 1"
$ws.Range("I4").Value = "Shift(6-14)"
$ws.Range("J4").Value = "Admin User"
$ws.Range("K4").Value = "--"
$ws.Range("L4").Value = "--"
$ws.Range("M4").Value = ""
$ws.Range("O4").Value = "EN MAINTENANCE"
